$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync
#    with en-US" everywhere it appears (Overview + per-locale sheets).
#  - The previously-empty "Latest Handback DateTime" placeholder
#    ("0001-01-01 00:00:00") is replaced with real handback timestamps.
#  - The "Latest Target File" / "Latest Handback File" columns (F/G) get
#    populated with hyperlinks pointing at the handed-back files.
# ---------------------------------------------------------------------------

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdDisplay  = "d828a29d-1351-4e93-ad4f-5f647b07e4aa.md"
$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/ac66369312e9e2e466043967109e6701fe321018/e2e/d828a29d-1351-4e93-ad4f-5f647b07e4aa.md"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("H2").Value = "2016-03-19 14:51:44"
$wsZh.Range("H3").Value = "2016-03-19 14:51:44"

$zhXlfDisplay = "d828a29d-1351-4e93-ad4f-5f647b07e4aa.473b3deb737bdc7f1ef6f31f672985aa1c11fd8c.zh-cn.xlf"
$zhXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/147fcce13af86c7861ac39ad7e5e90242b12a64b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d828a29d-1351-4e93-ad4f-5f647b07e4aa.473b3deb737bdc7f1ef6f31f672985aa1c11fd8c.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrl, "", "", $mdDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", $zhXlfDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrl, "", "", $mdDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", $zhXlfDisplay) | Out-Null

# --- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("H2").Value = "2016-03-19 14:51:50"
$wsDe.Range("H3").Value = "2016-03-19 14:51:50"

$deXlfDisplay = "d828a29d-1351-4e93-ad4f-5f647b07e4aa.473b3deb737bdc7f1ef6f31f672985aa1c11fd8c.de-de.xlf"
$deXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/197b252a493d0b9edbb3de43e3c422fdfe064bcd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d828a29d-1351-4e93-ad4f-5f647b07e4aa.473b3deb737bdc7f1ef6f31f672985aa1c11fd8c.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrl, "", "", $mdDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", $deXlfDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrl, "", "", $mdDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", $deXlfDisplay) | Out-Null
